$wb = $excel.ActiveWorkbook

# --- Rename headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy formatting (bold header / borders, date format for column A) from the
# "Weekly Quantity" sheet so the new sheet matches the look of the others.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Data rows ---
$wsForecast.Range("A2").Value = 45004.99999999999
$wsForecast.Range("B2").Value = 50
$wsForecast.Range("C2").Value = 49.99980579552651
$wsForecast.Range("D2").Value = 49.99980581293599

$wsForecast.Range("A3").Value = 45011.99999999999
$wsForecast.Range("B3").Value = 30
$wsForecast.Range("C3").Value = 29.99980579154797
$wsForecast.Range("D3").Value = 29.9998058084728

$wsForecast.Range("A4").Value = 45018.99999999999
$wsForecast.Range("B4").Value = 10
$wsForecast.Range("C4").Value = 9.999805398024792
$wsForecast.Range("D4").Value = 9.999806230282628

$wsForecast.Range("A5").Value = 45025.99999999999
$wsForecast.Range("B5").Value = 0
$wsForecast.Range("C5").Value = -10.00019546644152
$wsForecast.Range("D5").Value = -10.00019279596025

$wsForecast.Range("A6").Value = 45032.99999999999
$wsForecast.Range("B6").Value = 0
$wsForecast.Range("C6").Value = -30.00019658195261
$wsForecast.Range("D6").Value = -30.00019154410304

$wsForecast.Range("A7").Value = 45039.99999999999
$wsForecast.Range("B7").Value = 0
$wsForecast.Range("C7").Value = -50.00019774918807
$wsForecast.Range("D7").Value = -50.00019025680801

$wsForecast.Range("A8").Value = 45046.99999999999
$wsForecast.Range("B8").Value = 0
$wsForecast.Range("C8").Value = -70.00019921213342
$wsForecast.Range("D8").Value = -70.00018870969937

$wsForecast.Range("A9").Value = 45053.99999999999
$wsForecast.Range("B9").Value = 0
$wsForecast.Range("C9").Value = -90.0002008925288
$wsForecast.Range("D9").Value = -90.00018695071725

$wsForecast.Range("A10").Value = 45060.99999999999
$wsForecast.Range("B10").Value = 0
$wsForecast.Range("C10").Value = -110.000202668802
$wsForecast.Range("D10").Value = -110.0001851750264

$wsForecast.Range("A11").Value = 45067.99999999999
$wsForecast.Range("B11").Value = 0
$wsForecast.Range("C11").Value = -130.0002048215606
$wsForecast.Range("D11").Value = -130.0001827382752

$wsForecast.Range("A1").Select() | Out-Null
